$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update description text, widen row 3, adjust selection ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "Persons that enter the follow up in the cohort of {ImmDis}, with all their flares happened after entrance in the cohort (including before starting followup)"
$wsMeta.Range("B6").Value = "as many as the flares the person is experiencing after entering the cohort"
$wsMeta.Rows.Item(3).RowHeight = 28.8

# --- Data Model sheet: update the Rule text for the flare date row ---
$wsData = $wb.Worksheets.Item("Data Model")
$wsData.Range("K3").Value = "retrieve all distinct pairs (person_id date_flare) from D3_components_flare_TD_{ImmDis}, after removing record whose date_flare is missing"

# --- Selections / active sheet: make "Data Model" the active tab with K3 selected ---
$wsData.Activate()
$wsData.Range("K3").Select()

# --- Metadata sheet selection: A2:B13 range selected with A2 as the active cell ---
$wsMeta.Activate()
$wsMeta.Range("A2:B13").Select()

# --- Re-activate "Data Model" so it ends up as the active tab (activeTab=1) ---
$wsData.Activate()
